# Rewrites the first paragraph ("2" + bookmarked "1") into a clean
# "2"+"1" paragraph (no bookmark), then appends seven new paragraphs:
# three plain-text paragraphs and four "misspelled word" paragraphs
# wrapped in <w:proofErr> spell-check markers, the third of which
# carries the relocated "_GoBack" bookmark, matching the target diff.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Rebuild paragraph 1: drop the _GoBack bookmark and make the
#     trailing "1" a plain run (no rsid) right after "2". ---
$p1Xml = '<w:p ' + $wNs + '>' +
           '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>2</w:t></w:r>' +
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>1</w:t></w:r>' +
         '</w:p>'
$d.Paragraphs(1).Range.InsertXML($p1Xml)

# --- Helper-style construction of the new paragraphs' XML (plain text). ---
function Plain-ParaXml($text) {
    return '<w:p ' + $wNs + '>' +
             '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
             '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>' + $text + '</w:t></w:r>' +
           '</w:p>'
}

# --- Helper-style construction for "misspelled" paragraphs wrapped in
#     proofErr spellStart/spellEnd, optionally carrying the _GoBack bookmark. ---
function Spell-ParaXml($text, [bool]$withBookmark) {
    $bm = ''
    if ($withBookmark) {
        $bm = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
    }
    return '<w:p ' + $wNs + '>' +
             '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
             '<w:proofErr w:type="spellStart"/>' +
             '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>' + $text + '</w:t></w:r>' +
             $bm +
             '<w:proofErr w:type="spellEnd"/>' +
           '</w:p>'
}

# --- Append the seven new paragraphs, one at a time, after the current
#     last paragraph in the document. ---
$newParaXml = @(
    (Plain-ParaXml '12345678'),
    (Plain-ParaXml '12345678'),
    (Plain-ParaXml '123456789'),
    (Spell-ParaXml 'Sdfghgdfsasfdgvb' $false),
    (Spell-ParaXml 'Asfzadgfhgjhkjkl' $false),
    (Spell-ParaXml 'Asfdgfhgjhvjk' $true),
    (Spell-ParaXml 'zsdcgjhjkjvhjcgh' $false)
)

foreach ($xml in $newParaXml) {
    $endRange = $d.Paragraphs.Last.Range
    $endRange.Collapse(0)
    $endRange.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last.Range
    $newPara.InsertXML($xml)
}
